$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The size projection "aura" of the Team Agent ..." -> "The size of the
#    projection "aura" ..." plus a trailing period and two trailing spaces.
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "The size projection “aura” of the Team Agent will inversely affect the “rate of energy burn”",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The size of the projection “aura” of the Team Agent will inversely affect the “rate of energy burn”.  ",
    2)
if (-not $found1) { throw "change 1: source text not found" }

# ---------------------------------------------------------------------------
# 2) "The size of the detection"aura" projected by the ghosts ..." paragraph:
#    fix wording ("detection" -> "projection"), add two trailing sentences,
#    then append two brand-new bullet paragraphs (same list, numId=2) after
#    it: the "detector aura" definition and the "auras touch" note.
# ---------------------------------------------------------------------------
$old2 = "The size of the detection“aura” projected by the ghosts is a factor in the hardness of the game. The “easy” game will have a projection aura of at least four hallway widths for each ghost. The “easy” Basic Agent example will have about half the size of the projection aura of the ghosts. "
$new2 = "The size of the projection “aura” projected by the ghosts is a factor in the hardness of the game. The “easy” game will have a projection aura of at least four hallway widths for each ghost. The “easy” Basic Agent example will have about half the size of the projection aura of the ghosts. Projection auras are circular. They project thru walls. "

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) { throw "change 2: source text not found" }

# Locate the paragraph we just rewrote so we can append two new bullets after it.
$auraParaIdx = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Projection auras are circular*") {
        $auraParaIdx = $i
        break
    }
}
if ($auraParaIdx -eq 0) { throw "change 2: rewritten paragraph not found" }

$auraPara = $d.Paragraphs.Item($auraParaIdx)
$auraPara.Range.InsertParagraphAfter()
$detectorPara = $d.Paragraphs.Item($auraParaIdx + 1)
$detectorPara.Range.Text = "The “detector aura” (ie, how far Agents and ghosts see) is an ellipse with the Agent in one vertex and the other is the direction the Agent is looking. "

$detectorPara.Range.InsertParagraphAfter()
$touchPara = $d.Paragraphs.Item($auraParaIdx + 2)
$touchPara.Range.Text = "When “projection” and “detector” auras touch, AI’s are influenced. Think hit box overlap."

# ---------------------------------------------------------------------------
# 3) The blank "  " bullet (numId=1, just under "weigh and measure") becomes
#    the "collide in the finals" question, followed by a brand-new empty
#    bullet paragraph (numId=1).
# ---------------------------------------------------------------------------
$blankParaIdx = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "  `r") {
        $blankParaIdx = $i
        break
    }
}
if ($blankParaIdx -eq 0) { throw "change 3: blank bullet paragraph not found" }

$blankPara = $d.Paragraphs.Item($blankParaIdx)
$blankPara.Range.Text = "What happens when two Agents collide in the finals?  What happens when three or more collide at the same time, or slightly different times? Are two Agents able to gang up on a third Agent in the finals? "
$blankPara.Range.InsertParagraphAfter()

Write-Output "found1=$found1 found2=$found2 auraParaIdx=$auraParaIdx blankParaIdx=$blankParaIdx"
